# LogicComponentSequenceDiagram.pptx update
# Reflects renaming Address/BookParser -> ModulePlanner/Parser and
# delete-person-by-index -> delete-modules-by-course-code throughout
# the sequence diagram on slide 1.

$ldq = [char]0x201C   # left curly quote "
$rdq = [char]0x201D   # right curly quote "

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------
# Shape id=16 "Rectangle 62" : ":Address" / "BookParser"
#                            -> ":ModulePlanner" / "Parser"
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(5)
$tr = $shp.TextFrame.TextRange
$bookParser = $tr.Characters(10, 10)
$bookParser.Text = "Parser"
$address = $tr.Characters(2, 7)
$address.Text = "ModulePlanner"
$shp.Left = 203.999969482421875
$shp.Top = 171.5991058349609375
$shp.Width = 119.24224853515625
$shp.Height = 36.825473785400390625

# ---------------------------------------------------------------
# Shape id=25 "Straight Arrow Connector 24": reposition, drop vertical flip
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(12)
$shp.VerticalFlip = 0
$shp.Left = 128.6298828125
$shp.Top = 245.7068939208984375
$shp.Width = 130.0670623779296875
$shp.Height = 0.00007870377885410562

# ---------------------------------------------------------------
# Shape id=26 "TextBox 25": execute("delete 1") -> execute("delete / c/CS1010")
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(13)
$tr = $shp.TextFrame.TextRange
$tr.Text = "execute($ldq" + "delete`rc/CS1010$rdq)"
$tr.ParagraphFormat.Alignment = 2   # ppAlignCenter
$shp.Left = 7.80736255645751953125
$shp.Top = 245.8343048095703125
$shp.Width = 112.19256591796875
$shp.Height = 33.928073883056640625

# ---------------------------------------------------------------
# Shape id=28 "Straight Arrow Connector 27": reposition, drop vertical flip
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(14)
$shp.VerticalFlip = 0
$shp.Left = 413.95416259765625
$shp.Top = 291.591796875
$shp.Width = 75.93280029296875
$shp.Height = 0.22523622214794158936

# ---------------------------------------------------------------
# Shape id=78 "TextBox 77": deletePerson(p) -> deleteModules(m)
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(27)
$tr = $shp.TextFrame.TextRange
$argRun = $tr.Characters(13, 3)
$argRun.Text = "(m)"
$nameRun = $tr.Characters(1, 12)
$nameRun.Text = "deleteModules"

# ---------------------------------------------------------------
# Shape id=79 "TextBox 78": parse("1") -> parse("c/CS1010")
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(28)
$tr = $shp.TextFrame.TextRange
$tr.Text = "parse($ldq" + "c/CS1010$rdq)"
$shp.Left = 290.60516357421875
$shp.Top = 287.999969482421875
$shp.Width = 93.39476776123046875
$shp.Height = 14.54059123992919921875

# ---------------------------------------------------------------
# Shape id=80 "TextBox 79": parse("delete 1") -> parse("delete / c/CS1010")
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(29)
$tr = $shp.TextFrame.TextRange
$tr.Text = "parse($ldq" + "delete`rc/CS1010$rdq)"
$tr.ParagraphFormat.Alignment = 2   # ppAlignCenter
$shp.Left = 138.3324127197265625
$shp.Top = 250.2395782470703125
$shp.Width = 112.19256591796875
$shp.Height = 33.928073883056640625

# ---------------------------------------------------------------
# Shape id=51 "Straight Arrow Connector 50": reposition
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(41)
$shp.Left = 272.403045654296875
$shp.Top = 329.533294677734375
$shp.Width = 134.4158782958984375
$shp.Height = 0.0
